$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "2026-01-23T07:13:56.034Z"
$ws.Range("B14").Value = "John Doe"
$ws.Range("C14").Value = "Male"
$ws.Range("D14").Value = 36
$ws.Range("E14").Value = "Asthma Attack"
$ws.Range("F14").Value = "68e36fbc9950693f31349cc7"
$ws.Range("G14").Value = "John Doe is currently experiencing an active asthma attack with multiple severe symptoms affecting breathing and daily activities. Vital signs are stable at the moment, but immediate attention is advised due to symptom severity."
$ws.Range("H14").Value = "HIGH: Severe symptoms such as faster breathing, inability to perform usual tasks, and inability to sleep indicate urgent need for medical review."
$ws.Range("I14").Value = "Blood pressure is steady at 120/80 mmHg. | Heart rate remains normal at 70 bpm. | Oxygen saturation is within normal limits at 98%. | Body temperature is normal at 98.6°F. | Blood sugar levels are stable at 100 mg/dL."
$ws.Range("J14").Value = "Routine activities such as breakfast are maintained but symptom severity impacts the ability to perform usual tasks. | Increased coughing and difficulty breathing worsen during inactive periods, affecting sleep quality. | Symptoms like faster breathing significantly disrupt daily comfort."
$ws.Range("K14").Value = "Avoid exposure to asthma triggers like smoke, dust, and strong odors. | Maintain a calm environment to help ease breathing difficulties. | Ensure regular hydration and balanced nutrition to support overall health. | Keep track of symptoms and seek help promptly if they worsen."
$ws.Range("L14").Value = "Patient shows severe respiratory symptoms indicating an active asthma episode. | Stable vitals are encouraging but close monitoring is essential. | Recommend urgent consultation with a healthcare provider."
$ws.Range("M14").Value = "Schedule an immediate follow-up appointment with your asthma specialist. | Consider visiting an emergency facility if symptoms such as inability to talk or eat worsen. | Continue monitoring vitals and symptom severity closely at home."
$ws.Range("N14").Value = 1048
$ws.Range("O14").Value = 386
$ws.Range("P14").Value = 1434

$ws.Range("A15").Value = "2026-01-23T07:26:09.001Z"
$ws.Range("B15").Value = "John Doe"
$ws.Range("C15").Value = "Male"
$ws.Range("D15").Value = 36
$ws.Range("E15").Value = "Asthma Attack"
$ws.Range("F15").Value = "68e36fbc9950693f31349cc7"
$ws.Range("G15").Value = "John Doe is currently experiencing an active asthma attack with multiple severe symptoms impacting daily activities and sleep. Vital signs remain stable but symptoms suggest urgent attention is needed."
$ws.Range("H15").Value = "HIGH: Severe asthma symptoms such as inability to perform usual tasks, sleep, and eat require close monitoring and urgent medical attention. | MEDIUM: Increased coughing, difficulty breathing, and faster breathing indicate worsening asthma control."
$ws.Range("I15").Value = "Blood pressure, heart rate, oxygen level, body temperature, and blood sugar are within normal ranges. | Oxygen saturation is at 98%, which is normal but should be closely watched given respiratory symptoms."
$ws.Range("J15").Value = "Severe symptoms have significantly impacted the ability to perform daily tasks and sleep well. | Routine activities like meals may be affected due to inability to eat."
$ws.Range("K15").Value = "Avoid exposure to possible asthma triggers such as smoke or allergens. | Maintain a calm environment and practice controlled breathing exercises to ease breathing difficulty. | Ensure a balanced diet when possible to support overall health and immune function."
$ws.Range("L15").Value = "The patient is experiencing a severe asthma attack requiring immediate clinical review. | Symptom severity is high despite stable vital signs, indicating a need for urgent intervention."
$ws.Range("M15").Value = "Schedule an urgent consultation with a healthcare provider or asthma specialist. | Consider visiting the emergency room if symptoms worsen, especially inability to speak or eat, or worsening breathing difficulty. | Monitor symptoms closely and maintain regular vital sign checks. | Avoid known asthma triggers and prepare for possible lab tests or inhaler adjustments as directed by a healthcare professional."
$ws.Range("N15").Value = 1056
$ws.Range("O15").Value = 391
$ws.Range("P15").Value = 1447

$ws.Range("A16").Value = "2026-01-23T12:07:48.229Z"
$ws.Range("B16").Value = "John Doe"
$ws.Range("C16").Value = "Male"
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = "Asthma Attack"
$ws.Range("F16").Value = "68e36fbc9950693f31349cc7"
$ws.Range("G16").Value = "John Doe is currently experiencing an active asthma attack with significant respiratory symptoms impacting daily function. Vitals remain stable at present, but symptoms indicate severe respiratory distress requiring close monitoring."
$ws.Range("H16").Value = "HIGH: Severe breathing difficulties noted including inability to talk, eat, sleep, and perform usual tasks. Immediate follow-up is essential to prevent respiratory failure. | MEDIUM: Increased coughing, faster breathing, and gradual attack build-up suggest worsening airway obstruction."
$ws.Range("I16").Value = "Stable blood pressure at 120/80 mmHg. | Heart rate steady at 70 bpm. | Oxygen saturation normal at 98%. | Body temperature and blood sugar levels within normal limits."
$ws.Range("J16").Value = "Severe asthma symptoms are significantly limiting daily activities, including eating, sleeping, and speaking. | No variation noted in vitals throughout assessment period despite worsening symptoms."
$ws.Range("K16").Value = "Avoid known asthma triggers such as allergens, smoke, or strong odors. | Prioritize resting in an upright position to ease breathing. | Ensure good hydration and maintain a calm environment to reduce stress. | Seek immediate medical care if symptoms worsen or oxygen levels drop."
$ws.Range("L16").Value = "Patient is currently in respiratory distress with high symptom burden. | Vital signs do not yet reflect hypoxia but close monitoring is required. | Supportive care and prompt escalation plans need to be in place."
$ws.Range("M16").Value = "Urgent evaluation by respiratory specialist is recommended. | Continuous monitoring of oxygen saturation and respiratory status. | Consider referral to emergency care if symptoms fail to improve or worsen. | Educate patient on recognizing early warning signs of severe asthma exacerbations."
$ws.Range("N16").Value = 1138
$ws.Range("O16").Value = 407
$ws.Range("P16").Value = 1545
